$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.164.28"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").Value = "2.919.69"
$ws.Range("E3").Value = "  +3.43%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'204.32"
$ws.Range("E5").Value = "  +8.92%  "
$ws.Range("D6").Value = "'598.47"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  +2.70%  "
$ws.Range("D10").Value = "2.920.58"
$ws.Range("E10").Value = "  +3.50%  "
$ws.Range("E11").Value = "  +16.39%  "
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").Value = "3.456.59"
$ws.Range("E14").Value = "  +3.45%  "
$ws.Range("D15").Value = "76.076.89"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("D16").Value = "'28.10"
$ws.Range("E16").Value = "  +4.85%  "
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("D18").Value = "2.924.97"
$ws.Range("E18").Value = "  +3.85%  "
$ws.Range("D19").Value = "'12.90"
$ws.Range("E19").Value = "  +4.79%  "
$ws.Range("D20").Value = "'8.74"
$ws.Range("E20").Value = "  -2.23%  "
$ws.Range("D21").Value = "'372.81"
$ws.Range("E21").Value = "  -1.18%  "
$ws.Range("D22").Value = "'2.31"
$ws.Range("E22").Value = "  +2.22%  "
$ws.Range("D23").Value = "'4.29"
$ws.Range("E23").Value = "  +5.64%  "
$ws.Range("D24").Value = "'71.49"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "3.070.88"
$ws.Range("E26").Value = "  +3.56%  "
$ws.Range("E27").Value = "  +1.82%  "
$ws.Range("D28").Value = "'9.71"
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("E29").Value = "  +3.91%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("D32").Value = "'502.43"
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("D33").Value = "'7.75"
$ws.Range("E33").Value = "  +0.38%  "
$ws.Range("E34").Value = "  +2.46%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'165.21"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("D37").Value = "'20.21"
$ws.Range("E37").Value = "  +1.40%  "
$ws.Range("D38").Value = "'19.62"
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("E39").Value = "  +24.22%  "
$ws.Range("D40").Value = "'0.112"
$ws.Range("E40").Value = "  -4.84%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "'0.365"
$ws.Range("E41").Value = "  +6.94%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'182.79"
$ws.Range("E42").Value = "  -1.98%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "'5.01"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("D45").Value = "'1.66"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "'1.19"
$ws.Range("E47").Value = "  -1.48%  "
$ws.Range("D48").Value = "'2.38"
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("D49").Value = "'0.573"
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'22.39"
$ws.Range("E51").Value = "  +7.15%  "
